# Applies the "Finish product for sending" edit:
#  - LOAITIETKIEM: column F widened; F/G columns become descriptive-text flags
#    instead of raw 0/1 numbers (withdrawal rule + "Co" flag).
#  - SOTIETKIEM: account-number column F becomes text-formatted; one
#    obviously-wrong account number (row 3) corrected.
#  - THAMSO: a new leading index column is inserted and the old cryptic
#    parameter codes are replaced with human-readable Vietnamese labels.
#  - PHANQUYEN: selection anchor reset.
#  - THAMSO becomes the active sheet/tab for the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) LOAITIETKIEM (sheet "LOAITIETKIEM")
# ---------------------------------------------------------------------
$wsLoai = $wb.Worksheets.Item("LOAITIETKIEM")

# New column F width (stored width 12 == ColumnWidth 12 - 5/6)
$wsLoai.Columns.Item(6).ColumnWidth = 11.166666666666666

$wsLoai.Range("F1").Value = "Rút nhỏ hơn hoặc bằng"
$wsLoai.Range("G1").Value = "Có"

$wsLoai.Range("F2").Value = "Rút hết"
$wsLoai.Range("G2").Value = "Có"

$wsLoai.Range("F3").Value = "Rút hết"
$wsLoai.Range("G3").Value = "Có"

$wsLoai.Range("F1").Select()

# ---------------------------------------------------------------------
# 2) SOTIETKIEM (sheet "SOTIETKIEM")
# ---------------------------------------------------------------------
$wsSo = $wb.Worksheets.Item("SOTIETKIEM")

# Account numbers are long and should be stored/displayed as text.
$wsSo.Range("F1:F31").NumberFormat = "@"

# Row 3's account number had an extra bogus digit group - fix it.
$wsSo.Range("F3").Value = 15234634649

$wsSo.Range("A31").Select()

# ---------------------------------------------------------------------
# 3) THAMSO (sheet "THAMSO")
# ---------------------------------------------------------------------
$wsTham = $wb.Worksheets.Item("THAMSO")

# Insert a leading numbering column, pushing the existing two columns right.
$wsTham.Columns.Item(1).Insert()

$wsTham.Range("A1").Value = 1
$wsTham.Range("A2").Value = 2
$wsTham.Range("A3").Value = 3

$wsTham.Range("B1").Value = "Số tiền gửi tối thiểu"
$wsTham.Range("B2").Value = "Tiền gửi thêm tối thiểu"
$wsTham.Range("B3").Value = "Đóng sổ tự động"

$wsTham.Rows.Item(1).AutoFit()
$wsTham.Rows.Item(2).RowHeight = 30
$wsTham.Rows.Item(3).AutoFit()

# ---------------------------------------------------------------------
# 4) PHANQUYEN (sheet "PHANQUYEN")
# ---------------------------------------------------------------------
$wsPhan = $wb.Worksheets.Item("PHANQUYEN")
$wsPhan.Range("A1:C18").Select()

# ---------------------------------------------------------------------
# THAMSO ends up as the active sheet/tab - select it last.
# ---------------------------------------------------------------------
$wsTham.Activate()
$wsTham.Range("B3").Select()
